$wb = $excel.ActiveWorkbook

# Sheet "zh-cn": Priority (col E) low -> ht for rows 4-7,
# and Latest Handoff Datetime (col H) updated for rows 4-7 (shared timestamp)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("E5").Value = "ht"
$wsZh.Range("E6").Value = "ht"
$wsZh.Range("E7").Value = "ht"
$wsZh.Range("H4").Value = "2016-08-31 12:37:13"
$wsZh.Range("H5").Value = "2016-08-31 12:37:13"
$wsZh.Range("H6").Value = "2016-08-31 12:37:13"
$wsZh.Range("H7").Value = "2016-08-31 12:37:13"

# Sheet "de-de": Priority (col E) low -> ht for rows 4-7,
# and Latest Handoff Datetime (col H) updated for rows 4-7 (shared timestamp,
# same Xliff-generation run as the Overview sheet below)
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("E5").Value = "ht"
$wsDe.Range("E6").Value = "ht"
$wsDe.Range("E7").Value = "ht"
$wsDe.Range("H4").Value = "2016-08-31 12:37:18"
$wsDe.Range("H5").Value = "2016-08-31 12:37:18"
$wsDe.Range("H6").Value = "2016-08-31 12:37:18"
$wsDe.Range("H7").Value = "2016-08-31 12:37:18"

# Sheet "Overview": Latest HO Xliff Generate Date (col G) updated for rows 4-7 (shared timestamp)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-31 12:37:18"
$wsOverview.Range("G5").Value = "2016-08-31 12:37:18"
$wsOverview.Range("G6").Value = "2016-08-31 12:37:18"
$wsOverview.Range("G7").Value = "2016-08-31 12:37:18"
